$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column I (RF) for rows 25-53 from 34.21380952380952 to 20.389
$ws.Range("I25:I53").Value = 20.389
